$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column D (Price) as plain text so numeric-looking strings
# (e.g. "11.01", "0.0861") are not auto-coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "51.739.44"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "3.109.33"
$ws.Range("E3").Value = "  +4.15%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "388.68"
$ws.Range("E5").Value = "  +1.72%  "

$ws.Range("D6").Value = "103.35"
$ws.Range("E6").Value = "  -0.72%  "

$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "0.591"
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").Value = "37.36"
$ws.Range("E10").Value = "  +1.52%  "

$ws.Range("E11").Value = "  -0.09%  "

$ws.Range("D12").Value = "0.0861"

$ws.Range("D13").Value = "3.599.98"
$ws.Range("E13").Value = "  +3.89%  "

$ws.Range("D14").Value = "7.93"
$ws.Range("E14").Value = "  +1.37%  "

$ws.Range("D15").Value = "18.70"
$ws.Range("E15").Value = "  +1.08%  "

$ws.Range("D16").Value = "3.104.68"
$ws.Range("E16").Value = "  +3.94%  "

$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "11.01"
$ws.Range("E17").Value = "  -2.20%  "

$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "0.991"
$ws.Range("E18").Value = "  -0.33%  "

$ws.Range("D19").Value = "51.793.41"
$ws.Range("E19").Value = "  +0.33%  "

$ws.Range("E20").Value = "  +3.67%  "

$ws.Range("E21").Value = "  -0.82%  "

$ws.Range("E22").Value = "  +0.26%  "

$ws.Range("D23").Value = "70.04"
$ws.Range("E23").Value = "  -0.48%  "

$ws.Range("D24").Value = "267.68"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Value = "3.13"
$ws.Range("E25").Value = "  -3.15%  "

$ws.Range("D26").Value = "8.13"
$ws.Range("E26").Value = "  +1.53%  "

$ws.Range("E27").Value = "  +3.90%  "

$ws.Range("E28").Value = "  +0.38%  "

$ws.Range("E29").Value = "  +0.04%  "

$ws.Range("D30").Value = "7.15"
$ws.Range("E30").Value = "  -0.64%  "

$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").Value = "35.56"
$ws.Range("E33").Value = "  +2.57%  "

$ws.Range("D35").Value = "50.28"
$ws.Range("E35").Value = "  -2.12%  "

$ws.Range("D36").Value = "0.0449"
$ws.Range("E36").Value = "  +0.93%  "

$ws.Range("E37").Value = "  -0.13%  "

$ws.Range("E38").Value = "  +2.90%  "

$ws.Range("D39").Value = "0.291"
$ws.Range("E39").Value = "  +6.80%  "

$ws.Range("E40").Value = "  +2.43%  "

$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("D42").Value = "129.15"
$ws.Range("E42").Value = "  +5.16%  "

$ws.Range("E43").Value = "  -1.06%  "

$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("E45").Value = "  -3.51%  "

$ws.Range("D46").Value = "22.23"
$ws.Range("E46").Value = "  +3.76%  "

$ws.Range("E47").Value = "  +6.06%  "

$ws.Range("E48").Value = "  +2.26%  "

$ws.Range("D49").Value = "2.046.40"
$ws.Range("E49").Value = "  +0.56%  "

$ws.Range("D50").Value = "3.414.75"
$ws.Range("E50").Value = "  +4.02%  "

$ws.Range("D51").Value = "0.207"
$ws.Range("E51").Value = "  +5.09%  "
